$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "26.863.77"
$ws.Cells.Item(2, 5).Value = "  -0.89%  "

$ws.Cells.Item(3, 4).Value = "1.562.53"
$ws.Cells.Item(3, 5).Value = "  +0.04%  "

$ws.Cells.Item(4, 5).Value = "  -0.19%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "206.05"
$ws.Cells.Item(5, 5).Value = "  -0.11%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.488"
$ws.Cells.Item(6, 5).Value = "  -0.90%  "

$ws.Cells.Item(7, 5).Value = "  -0.22%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "21.72"
$ws.Cells.Item(8, 5).Value = "  -1.69%  "

$ws.Cells.Item(9, 5).Value = "  -0.04%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.0584"
$ws.Cells.Item(10, 5).Value = "  -1.24%  "

$ws.Cells.Item(11, 5).Value = "  +0.37%  "

$ws.Cells.Item(12, 4).Value = "1.784.86"
$ws.Cells.Item(12, 5).Value = "  +0.01%  "

$ws.Cells.Item(13, 4).Value = "1.565.19"

$ws.Cells.Item(14, 5).Value = "  -0.98%  "

$ws.Cells.Item(15, 5).Value = "  -0.13%  "

$ws.Cells.Item(16, 4).Value = "26.866.12"
$ws.Cells.Item(16, 5).Value = "  -0.92%  "

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "61.27"
$ws.Cells.Item(17, 5).Value = "  -2.70%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "214.59"
$ws.Cells.Item(18, 5).Value = "  +1.31%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "7.35"
$ws.Cells.Item(19, 5).Value = "  +2.10%  "

$ws.Cells.Item(20, 4).Value = "0.0₃0679"
$ws.Cells.Item(20, 5).Value = "  -1.08%  "

$ws.Cells.Item(21, 5).Value = "  -0.15%  "

$ws.Cells.Item(22, 5).Value = "  +0.43%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "9.16"
$ws.Cells.Item(23, 5).Value = "  -2.30%  "

$ws.Cells.Item(24, 5).Value = "  +1.23%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "153.91"
$ws.Cells.Item(25, 5).Value = "  +1.24%  "

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "6.74"
$ws.Cells.Item(26, 5).Value = "  +2.92%  "

$ws.Cells.Item(27, 5).Value = "  +0.43%  "

$ws.Cells.Item(28, 5).Value = "  -0.18%  "

$ws.Cells.Item(29, 5).Value = "  -0.93%  "

$ws.Cells.Item(30, 5).Value = "  +0.16%  "

$ws.Cells.Item(31, 5).Value = "  -3.40%  "

$ws.Cells.Item(32, 5).Value = "  +0.06%  "

$ws.Cells.Item(33, 4).Value = "1.402.65"
$ws.Cells.Item(33, 5).Value = "  +1.86%  "

$ws.Cells.Item(34, 5).Value = "  -0.29%  "

$ws.Cells.Item(35, 5).Value = "  -1.30%  "

$ws.Cells.Item(36, 5).Value = "  -0.39%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.917"
$ws.Cells.Item(37, 5).Value = "  -2.60%  "

$ws.Cells.Item(38, 5).Value = "  -0.29%  "

$ws.Cells.Item(39, 5).Value = "  +1.37%  "

$ws.Cells.Item(40, 5).Value = "  -0.27%  "

$ws.Cells.Item(41, 5).Value = "  -0.14%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.998"
$ws.Cells.Item(42, 5).Value = "  +0.64%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "5.39"
$ws.Cells.Item(43, 5).Value = "  +3.51%  "

$ws.Cells.Item(44, 5).Value = "  +0.15%  "

$ws.Cells.Item(45, 5).Value = "  -1.05%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "63.18"
$ws.Cells.Item(46, 5).Value = "  -0.28%  "

$ws.Cells.Item(47, 4).Value = "1.697.96"
$ws.Cells.Item(47, 5).Value = "  +0.12%  "

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "86.18"
$ws.Cells.Item(48, 5).Value = "  +1.02%  "

$ws.Cells.Item(49, 5).Value = "  +2.66%  "

$ws.Cells.Item(50, 4).Value = "0.0₇0976"
$ws.Cells.Item(50, 5).Value = "  -2.04%  "

$ws.Cells.Item(51, 5).Value = "  +0.72%  "
